# Damian Oguche.docx — "Updated the assets folder files"
#
# Three bullet/heading lines get their tech-stack / title wording expanded:
#   1. E-Commerce Platform bullet: call out React + Firebase + Stripe explicitly.
#   2. API Development bullet: call out Node, Express.js, TypeScript + MongoDB explicitly.
#   3. Polaris Bank Limited line: retitle "Team Lead, ATMs and POS Support Team"
#      to "Regional ATM/POS Support Officer".
#
# wdReplaceAll = 2, wdFindContinue (Wrap) = 1

$d = $word.ActiveDocument

function Replace-Exact($doc, $findText, $replaceText) {
    $rng = $doc.Content
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Could not find text to replace: $findText"
    }
    return $ok
}

# 1) E-Commerce Platform bullet.
Replace-Exact $d `
    " and React, integrating payment systems and real-time product updates." `
    ", React and Firebase, integrating Stripe payment systems and real-time product updates."

# 2) API Development bullet.
Replace-Exact $d `
    " Created RESTful APIs with Express.js for a logistics tracking system, enabling real-time updates for users." `
    " Created RESTful APIs with Node, Express.js, TypeScript and MongoDB for a logistics tracking system, enabling real-time updates for users."

# 3) Polaris Bank Limited role line.
Replace-Exact $d `
    " — Team Lead, ATMs and POS Support Team" `
    " — Regional ATM/POS Support Officer"

Write-Output "Done."
